$wb = $excel.ActiveWorkbook

# --- "About" sheet (sheet1) ---
$wsAbout = $wb.Worksheets.Item("About")

# Update the two headline labels (row 1-2)
$wsAbout.Range("A1").Value = "BpTPEU BTU per Large Primary Energy Unit"
$wsAbout.Range("A2").Value = "BpTPEU BTU per Small Primary Energy Unit"

# Update the unit description notes (row 11-12)
$wsAbout.Range("A11").Value = "The large primary energy output unit (used in totals graphs) is: quadrillion BTU"
$wsAbout.Range("A12").Value = "The small primary energy output unit (used in energy intensity per unit GDP graphs) is: thousand BTU"

# Remove the old "1 Btu = ..." footnote row entirely
$wsAbout.Range("A15").EntireRow.Delete()

# --- "BpTPEU-large" sheet (sheet2) ---
$wsLarge = $wb.Worksheets.Item("BpTPEU-large")
$wsLarge.Range("B2").Formula = "=10^15"

# --- "BpTPEU-small" sheet (sheet3) ---
$wsSmall = $wb.Worksheets.Item("BpTPEU-small")
$wsSmall.Range("B2").Formula = "=10^3"

# --- Active sheet / selection bookkeeping ---
# Reset selections on the two data sheets back to A1
$wsLarge.Range("A1").Select()
$wsSmall.Range("A1").Select()

# Make "About" the active sheet/tab again, with A1 selected
$wsAbout.Activate()
$wsAbout.Range("A1").Select()
